$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DESGLOSE")

$c1 = $ws.Cells.Item(5,1)
$c1.NumberFormat = "@"
$c1.Value = "100%"
$c1.NumberFormat = "General"
Write-Output "c1 value: $($c1.Value())"
